# Applies the scheduled market-data refresh described in the commit diff.
# Updates currentAveragePrice* / Leve*Price* / LeveProfit* columns (H:N)
# for the affected Leve rows across all eight sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 518.8
$ws.Range("I19").Value = 645
$ws.Range("J19").Value = 487.25
$ws.Range("K19").Value = 645
$ws.Range("L19").Value = 487.25
$ws.Range("M19").Value = -470
$ws.Range("N19").Value = -837.25

# Row 86
$ws.Range("H86").Value = 3203.6562
$ws.Range("I86").Value = 1428
$ws.Range("J86").Value = 5798.846
$ws.Range("K86").Value = 1428
$ws.Range("L86").Value = 5798.846
$ws.Range("M86").Value = -305
$ws.Range("N86").Value = -8044.846

# Row 89
$ws.Range("H89").Value = 3203.6562
$ws.Range("I89").Value = 1428
$ws.Range("J89").Value = 5798.846
$ws.Range("K89").Value = 7140
$ws.Range("L89").Value = 28994.23
$ws.Range("M89").Value = -1524
$ws.Range("N89").Value = -40226.23

# Row 132
$ws.Range("H132").Value = 2230.0127
$ws.Range("I132").Value = 1693.0807
$ws.Range("J132").Value = 4188.2354
$ws.Range("K132").Value = 5079.242099999999
$ws.Range("L132").Value = 12564.7062
$ws.Range("M132").Value = -2549.242099999999
$ws.Range("N132").Value = -17624.7062

# Row 138
$ws.Range("H138").Value = 2552.4727
$ws.Range("I138").Value = 1274.3334
$ws.Range("J138").Value = 4086.24
$ws.Range("K138").Value = 3823.0002
$ws.Range("L138").Value = 12258.72
$ws.Range("M138").Value = 1316.9998
$ws.Range("N138").Value = -22538.72

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 44069.566
$ws.Range("I97").Value = 63034.875
$ws.Range("K97").Value = 63034.875
$ws.Range("M97").Value = -62538.875

# Row 132
$ws.Range("H132").Value = 6454.825
$ws.Range("I132").Value = 4417.387
$ws.Range("J132").Value = 13472.667
$ws.Range("K132").Value = 13252.161
$ws.Range("L132").Value = 40418.001
$ws.Range("M132").Value = -10722.161
$ws.Range("N132").Value = -45478.001

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1270.6897
$ws.Range("I134").Value = 939.8889
$ws.Range("J134").Value = 1812
$ws.Range("K134").Value = 2819.6667
$ws.Range("L134").Value = 5436
$ws.Range("M134").Value = -284.6667000000002
$ws.Range("N134").Value = -10506

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1687.6957
$ws.Range("I16").Value = 739.4375
$ws.Range("J16").Value = 3855.1428
$ws.Range("K16").Value = 739.4375
$ws.Range("L16").Value = 3855.1428
$ws.Range("M16").Value = -452.4375
$ws.Range("N16").Value = -4429.1428

# Row 31
$ws.Range("H31").Value = 17280014
$ws.Range("I31").Value = 33334696
$ws.Range("J31").Value = 78567.75
$ws.Range("K31").Value = 33334696
$ws.Range("L31").Value = 78567.75
$ws.Range("M31").Value = -33334401
$ws.Range("N31").Value = -79157.75

# Row 34
$ws.Range("H34").Value = 17280014
$ws.Range("I34").Value = 33334696
$ws.Range("J34").Value = 78567.75
$ws.Range("K34").Value = 33334696
$ws.Range("L34").Value = 78567.75
$ws.Range("M34").Value = -33334494
$ws.Range("N34").Value = -78971.75

# Row 107
$ws.Range("H107").Value = 144174.58
$ws.Range("I107").Value = 334999.66
$ws.Range("J107").Value = 1055.75
$ws.Range("K107").Value = 334999.66
$ws.Range("L107").Value = 1055.75
$ws.Range("M107").Value = -333079.66
$ws.Range("N107").Value = -4895.75

# Row 113
$ws.Range("H113").Value = 1687.6957
$ws.Range("I113").Value = 739.4375
$ws.Range("J113").Value = 3855.1428
$ws.Range("K113").Value = 739.4375
$ws.Range("L113").Value = 3855.1428
$ws.Range("M113").Value = 1430.5625
$ws.Range("N113").Value = -8195.1428

# Row 132
$ws.Range("H132").Value = 15876932
$ws.Range("I132").Value = 22731174
$ws.Range("J132").Value = 3952.6316
$ws.Range("K132").Value = 68193522
$ws.Range("L132").Value = 11857.8948
$ws.Range("M132").Value = -68190992
$ws.Range("N132").Value = -16917.8948

$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 75
$ws.Range("H75").Value = 870.25
$ws.Range("I75").Value = 354.33334
$ws.Range("J75").Value = 1042.2222
$ws.Range("K75").Value = 1063.00002
$ws.Range("L75").Value = 3126.6666
$ws.Range("M75").Value = -65.00001999999995
$ws.Range("N75").Value = -5122.6666

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 78
$ws.Range("H78").Value = 870.25
$ws.Range("I78").Value = 354.33334
$ws.Range("J78").Value = 1042.2222
$ws.Range("K78").Value = 3189.00006
$ws.Range("L78").Value = 9379.9998
$ws.Range("M78").Value = 1802.99994
$ws.Range("N78").Value = -19363.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2643.5679
$ws.Range("I132").Value = 2751.1667
$ws.Range("J132").Value = 2170.1333
$ws.Range("K132").Value = 8253.500100000001
$ws.Range("L132").Value = 6510.3999
$ws.Range("M132").Value = -5723.500100000001
$ws.Range("N132").Value = -11570.3999

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 25370.5
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 33660.668
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 33660.668
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -34250.668

# Row 27
$ws.Range("H27").Value = 25370.5
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 33660.668
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 33660.668
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -33874.668

# Row 93
$ws.Range("H93").Value = 2269.7058
$ws.Range("I93").Value = 1980.3846
$ws.Range("J93").Value = 3210
$ws.Range("K93").Value = 1980.3846
$ws.Range("L93").Value = 3210
$ws.Range("M93").Value = -732.3846000000001
$ws.Range("N93").Value = -5706

# Row 132
$ws.Range("H132").Value = 3958.3174
$ws.Range("I132").Value = 4523.4736
$ws.Range("J132").Value = 3099.28
$ws.Range("K132").Value = 13570.4208
$ws.Range("L132").Value = 9297.84
$ws.Range("M132").Value = -11040.4208
$ws.Range("N132").Value = -14357.84

# Row 136
$ws.Range("H136").Value = 2533.4084
$ws.Range("I136").Value = 1050.6346
$ws.Range("J136").Value = 6591.5264
$ws.Range("K136").Value = 3151.9038
$ws.Range("L136").Value = 19774.5792
$ws.Range("M136").Value = -601.9038
$ws.Range("N136").Value = -24874.5792

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 10641916
$ws.Range("I132").Value = 16133444
$ws.Range("J132").Value = 2078.125
$ws.Range("K132").Value = 48400332
$ws.Range("L132").Value = 6234.375
$ws.Range("M132").Value = -48397802
$ws.Range("N132").Value = -11294.375

# Row 136
$ws.Range("H136").Value = 982.93445
$ws.Range("I136").Value = 459.51163
$ws.Range("J136").Value = 2233.3333
$ws.Range("K136").Value = 1378.53489
$ws.Range("L136").Value = 6699.999899999999
$ws.Range("M136").Value = 1171.46511
$ws.Range("N136").Value = -11799.9999
